$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("site paulo")

# Clear the two text cells first so the shared-string table rebuilds
# its entries for B3/C3 in the same order as the target workbook
# (duration text first, then the longer description text).
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""

# C3 now holds the short duration ("2h"), B3 holds the updated
# description of the task.
$ws.Range("C3").Value = "2h"
$ws.Range("B3").Value = "Criação de E-mail em HTML - versão 1 e alterações para 2."

# D3 now shows 2 decimal hours worked (previously 1.25)
$ws.Range("D3").Value = 2

# Move the active selection to B5 (as recorded in the saved file)
$ws.Range("B5").Select()

$wb.Application.Calculate()
